$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text (avoids Excel auto-converting
# numeric-looking strings like "315.00" or "1.002" into real numbers),
# then restore the default "Normal" style so no stray style index is left
# behind on the cell (matches original unstyled cells).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Row-by-row price / volume updates ---
Set-TextValue 'D2' '28.421.59'
Set-TextValue 'E2' '  +0.20%  '
Set-TextValue 'D3' '1.817.60'
Set-TextValue 'E3' '  -0.35%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '315.00'
Set-TextValue 'E5' '  -0.80%  '
Set-TextValue 'D6' '1.002'
Set-TextValue 'E6' '  +0.10%  '
Set-TextValue 'D7' '0.5107'
Set-TextValue 'E7' '  -4.35%  '
Set-TextValue 'D8' '0.3952'
Set-TextValue 'E8' '  -2.08%  '
Set-TextValue 'D9' '0.08107'
Set-TextValue 'E9' '  +6.46%  '
Set-TextValue 'E10' '  -0.33%  '
Set-TextValue 'D11' '1.107'
Set-TextValue 'E11' '  -0.08%  '
Set-TextValue 'D12' '20.97'
Set-TextValue 'E12' '  +0.35%  '
Set-TextValue 'D13' '6.284'
Set-TextValue 'E13' '  -0.56%  '
Set-TextValue 'D14' '1.002'
Set-TextValue 'E14' '  +0.12%  '
Set-TextValue 'D15' '7.507'
Set-TextValue 'E15' '  -1.30%  '
Set-TextValue 'D16' '1.818.34'
Set-TextValue 'E16' '  -0.62%  '
Set-TextValue 'D17' '0.00001133'
Set-TextValue 'E17' '  +5.34%  '
Set-TextValue 'D18' '92.60'
Set-TextValue 'E18' '  +3.58%  '
Set-TextValue 'E19' '  +0.62%  '
Set-TextValue 'E20' '  +0.12%  '
Set-TextValue 'E21' '  +0.06%  '
Set-TextValue 'D22' '6.098'
Set-TextValue 'E22' '  +0.11%  '
Set-TextValue 'D23' '28.456.70'
Set-TextValue 'E23' '  +0.26%  '
Set-TextValue 'D24' '11.27'
Set-TextValue 'E24' '  +0.89%  '
Set-TextValue 'D25' '2.269'
Set-TextValue 'E25' '  +2.86%  '
Set-TextValue 'D26' '21.16'
Set-TextValue 'E26' '  +2.61%  '
Set-TextValue 'D29' '2.403'
Set-TextValue 'E29' '  -2.24%  '
Set-TextValue 'E30' '  +1.68%  '
Set-TextValue 'E31' '  +0.03%  '
Set-TextValue 'D32' '1.105'
Set-TextValue 'E32' '  -1.42%  '
Set-TextValue 'D33' '5.788'
Set-TextValue 'E33' '  +2.31%  '
Set-TextValue 'D34' '3.652'
Set-TextValue 'E34' '  +0.25%  '
Set-TextValue 'D35' '0.07009'
Set-TextValue 'E35' '  -5.68%  '
Set-TextValue 'D36' '0.2225'
Set-TextValue 'E36' '  -0.37%  '
Set-TextValue 'D37' '5.232'
Set-TextValue 'E37' '  +0.66%  '
Set-TextValue 'D38' '0.02330'
Set-TextValue 'E38' '  -0.52%  '
Set-TextValue 'D39' '8.817'
Set-TextValue 'E39' '  -1.01%  '
Set-TextValue 'D40' '0.6280'
Set-TextValue 'E40' '  +0.40%  '
Set-TextValue 'D41' '11.29'
Set-TextValue 'E41' '  -0.12%  '
Set-TextValue 'D42' '1.174'
Set-TextValue 'E42' '  -0.62%  '
Set-TextValue 'D43' '1.001'
Set-TextValue 'E43' '  +0.12%  '
Set-TextValue 'D44' '1.402'
Set-TextValue 'E44' '  +0.54%  '
Set-TextValue 'D45' '13.51'
Set-TextValue 'E45' '  -0.10%  '
Set-TextValue 'D46' '3.744'
Set-TextValue 'E46' '  +1.19%  '
Set-TextValue 'D47' '0.5924'
Set-TextValue 'E47' '  +1.42%  '
Set-TextValue 'D48' '124.81'
Set-TextValue 'E48' '  -0.18%  '
Set-TextValue 'D49' '1.974'
Set-TextValue 'E49' '  -0.75%  '
Set-TextValue 'D50' '1.188'
Set-TextValue 'E50' '  -0.93%  '
Set-TextValue 'D51' '0.06886'
Set-TextValue 'E51' '  -0.14%  '

# --- Rows 27/28: Monero and WrappedliquidstakedEther2.0 swap positions ---
Set-TextValue 'B27' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C27' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D27' '2.036.15'
Set-TextValue 'E27' '  -0.22%  '
Set-TextValue 'B28' 'Monero'
Set-TextValue 'C28' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D28' '155.52'
Set-TextValue 'E28' '  -1.24%  '
